$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3 ("Détails sponsor") with new values matching the latest commit
$ws.Range("A3").Value = "sdfsd"
$ws.Range("B3").Value = "sdf"
$ws.Range("C3").Value = "sfdfds"
$ws.Range("D3").Value = 4774.0
$ws.Range("E3").Value = "dsds"
